$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9987.333000000001
$ws.Range("I51").Value = 8495
$ws.Range("J51").Value = 10173.875
$ws.Range("K51").Value = 8495
$ws.Range("L51").Value = 10173.875
$ws.Range("M51").Value = -8011
$ws.Range("N51").Value = -11141.875

$ws.Range("H70").Value = 11076.421
$ws.Range("I70").Value = 2432.875
$ws.Range("J70").Value = 17362.637
$ws.Range("K70").Value = 7298.625
$ws.Range("L70").Value = 52087.91099999999
$ws.Range("M70").Value = -7028.625
$ws.Range("N70").Value = -52627.91099999999

$ws.Range("H73").Value = 11076.421
$ws.Range("I73").Value = 2432.875
$ws.Range("J73").Value = 17362.637
$ws.Range("K73").Value = 7298.625
$ws.Range("L73").Value = 52087.91099999999
$ws.Range("M73").Value = -6362.625
$ws.Range("N73").Value = -53959.91099999999

$ws.Range("H92").Value = 2646.611
$ws.Range("I92").Value = 3035.5715
$ws.Range("K92").Value = 3035.5715
$ws.Range("M92").Value = -1787.5715

$ws.Range("H116").Value = 4942.5
$ws.Range("I116").Value = 4635.1
$ws.Range("K116").Value = 4635.1
$ws.Range("M116").Value = -1193.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4387.0513
$ws.Range("I32").Value = 1909.3485
$ws.Range("J32").Value = 18014.416
$ws.Range("K32").Value = 1909.3485
$ws.Range("L32").Value = 18014.416
$ws.Range("M32").Value = -1622.3485
$ws.Range("N32").Value = -18588.416

$ws.Range("H38").Value = 12217.857
$ws.Range("I38").Value = 19999
$ws.Range("J38").Value = 10921
$ws.Range("K38").Value = 19999
$ws.Range("L38").Value = 10921
$ws.Range("M38").Value = -19532
$ws.Range("N38").Value = -11855

$ws.Range("H61").Value = 3832.5417
$ws.Range("I61").Value = 2792.7222
$ws.Range("K61").Value = 2792.7222
$ws.Range("M61").Value = -2580.7222

$ws.Range("H122").Value = 2699
$ws.Range("I122").Value = 1563.4286
$ws.Range("K122").Value = 4690.2858
$ws.Range("M122").Value = -2240.2858

$ws.Range("H136").Value = 3832.5417
$ws.Range("I136").Value = 2792.7222
$ws.Range("K136").Value = 8378.1666
$ws.Range("M136").Value = -5828.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 803.15625
$ws.Range("I80").Value = 1090.875
$ws.Range("J80").Value = 515.4375
$ws.Range("K80").Value = 1090.875
$ws.Range("L80").Value = 515.4375
$ws.Range("M80").Value = -92.875
$ws.Range("N80").Value = -2511.4375

$ws.Range("H83").Value = 803.15625
$ws.Range("I83").Value = 1090.875
$ws.Range("J83").Value = 515.4375
$ws.Range("K83").Value = 5454.375
$ws.Range("L83").Value = 2577.1875
$ws.Range("M83").Value = -462.375
$ws.Range("N83").Value = -12561.1875

$ws.Range("H107").Value = 4513.1562
$ws.Range("I107").Value = 4610.387
$ws.Range("J107").Value = 1499
$ws.Range("K107").Value = 4610.387
$ws.Range("L107").Value = 1499
$ws.Range("M107").Value = -2690.387
$ws.Range("N107").Value = -5339

$ws.Range("H134").Value = 4566.324
$ws.Range("I134").Value = 4589.1875
$ws.Range("K134").Value = 13767.5625
$ws.Range("M134").Value = -11232.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2798.5715
$ws.Range("I58").Value = 2032.6957
$ws.Range("K58").Value = 2032.6957
$ws.Range("M58").Value = -1829.6957

$ws.Range("H62").Value = 3861
$ws.Range("I62").Value = 3861
$ws.Range("K62").Value = 3861
$ws.Range("M62").Value = -3237

$ws.Range("H65").Value = 3861
$ws.Range("I65").Value = 3861
$ws.Range("K65").Value = 19305
$ws.Range("M65").Value = -16185

$ws.Range("H86").Value = 1434200.2
$ws.Range("I86").Value = 2005780.8
$ws.Range("K86").Value = 2005780.8
$ws.Range("M86").Value = -2004657.8

$ws.Range("H89").Value = 1434200.2
$ws.Range("I89").Value = 2005780.8
$ws.Range("K89").Value = 10028904
$ws.Range("M89").Value = -10023288

$ws.Range("H99").Value = 427034.03
$ws.Range("I99").Value = 914447.2
$ws.Range("K99").Value = 914447.2
$ws.Range("M99").Value = -912949.2

$ws.Range("H122").Value = 2801.5
$ws.Range("I122").Value = 2993.1
$ws.Range("K122").Value = 8979.299999999999
$ws.Range("M122").Value = -6529.299999999999

$ws.Range("H126").Value = 427034.03
$ws.Range("I126").Value = 914447.2
$ws.Range("K126").Value = 2743341.6
$ws.Range("M126").Value = -2740871.6

$ws.Range("H134").Value = 1924.3654
$ws.Range("I134").Value = 1924.3654
$ws.Range("K134").Value = 5773.0962
$ws.Range("M134").Value = -3238.0962

$ws.Range("H136").Value = 2798.5715
$ws.Range("I136").Value = 2032.6957
$ws.Range("K136").Value = 6098.0871
$ws.Range("M136").Value = -3548.0871

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2607784
$ws.Range("J68").Value = 3962.5356
$ws.Range("L68").Value = 11887.6068
$ws.Range("N68").Value = -13509.6068

$ws.Range("H71").Value = 2607784
$ws.Range("J71").Value = 3962.5356
$ws.Range("L71").Value = 35662.8204
$ws.Range("N71").Value = -43774.8204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8000
$ws.Range("J70").Value = 8000
$ws.Range("L70").Value = 8000
$ws.Range("N70").Value = -8540

$ws.Range("H73").Value = 8000
$ws.Range("J73").Value = 8000
$ws.Range("L73").Value = 8000
$ws.Range("N73").Value = -9872

$ws.Range("I113").Value = 5186.9
$ws.Range("J113").Value = 5968.5
$ws.Range("K113").Value = 5186.9
$ws.Range("L113").Value = 5968.5
$ws.Range("M113").Value = -3016.9
$ws.Range("N113").Value = -10308.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9162.444
$ws.Range("I40").Value = 10247.172
$ws.Range("J40").Value = 5365.9
$ws.Range("K40").Value = 10247.172
$ws.Range("L40").Value = 5365.9
$ws.Range("M40").Value = -10111.172
$ws.Range("N40").Value = -5637.9

$ws.Range("H46").Value = 3884.4666
$ws.Range("I46").Value = 3584.5
$ws.Range("K46").Value = 3584.5
$ws.Range("M46").Value = -3396.5

$ws.Range("H122").Value = 4650
$ws.Range("I122").Value = 4650
$ws.Range("K122").Value = 13950
$ws.Range("M122").Value = -11500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 5378.2
$ws.Range("J45").Value = 5473.5
$ws.Range("L45").Value = 5473.5
$ws.Range("N45").Value = -6455.5

$ws.Range("H122").Value = 3781.5557
$ws.Range("I122").Value = 3754.3125
$ws.Range("K122").Value = 11262.9375
$ws.Range("M122").Value = -8812.9375

$ws.Range("H132").Value = 2154.1738
$ws.Range("J132").Value = 2638.889
$ws.Range("L132").Value = 7916.667
$ws.Range("N132").Value = -12976.667

$ws.Range("H136").Value = 2501501
$ws.Range("I136").Value = 5001002
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 15003006
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -15000456
$ws.Range("N136").Value = -11100
